$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 976
$ws1.Range("F6").Value = 1105
$ws1.Range("F7").Value = 885
$ws1.Range("F8").Value = 270
$ws1.Range("F11").Value = 870
$ws1.Range("F15").Value = 1358
$ws1.Range("F18").Value = 1233
$ws1.Range("F19").Value = 2913
$ws1.Range("F20").Value = 1501
$ws1.Range("F23").Value = 1296
$ws1.Range("F24").Value = 61
$ws1.Range("F26").Value = 368
$ws1.Range("F27").Value = 3227
$ws1.Range("F28").Value = 632
$ws1.Range("F29").Value = 544
$ws1.Range("F30").Value = 1440

# ---- Sheet "本地生活" (Local Life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 762

# ---- Sheet "全部类型" (All Types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 762
$ws4.Range("F6").Value = 976
$ws4.Range("F10").Value = 1105
$ws4.Range("F11").Value = 885
$ws4.Range("F12").Value = 270
$ws4.Range("F21").Value = 870
$ws4.Range("F25").Value = 1358
$ws4.Range("F28").Value = 1233
$ws4.Range("F29").Value = 2913
$ws4.Range("F30").Value = 1501
$ws4.Range("F33").Value = 1296
$ws4.Range("F34").Value = 61
$ws4.Range("F38").Value = 368
$ws4.Range("F39").Value = 3227
$ws4.Range("F40").Value = 632
$ws4.Range("F41").Value = 544
$ws4.Range("F42").Value = 1440
